$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.058311700820923
$ws.Range("B1").Value = 2.121825218200684
$ws.Range("C1").Value = 3.746902465820312
$ws.Range("D1").Value = 1.310323119163513
$ws.Range("E1").Value = 0.4232529699802399
